$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 1.4
$ws.Range("I3").Value = 9
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 8.5
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("X3").Value = 6
$ws.Range("AC3").Value = 8.5
$ws.Range("AE3").Value = 23
$ws.Range("AN3").Value = 3.2
$ws.Range("A4").Value = "OIegdTaq"
$ws.Range("C4").Value = "21:30"
$ws.Range("D4").Value = "ARGENTINA - TORNEO BETANO"
$ws.Range("E4").Value = "Sarmiento Junin"
$ws.Range("F4").Value = "Boca Juniors"
$ws.Range("G4").Value = 4.5
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 1.95
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 2.75
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 15
$ws.Range("Z4").Value = 51
$ws.Range("AA4").Value = 41
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 21
$ws.Range("AG4").Value = 5.5
$ws.Range("AH4").Value = 8
$ws.Range("AI4").Value = 9.5
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 19
$ws.Range("AM4").Value = 501
$ws.Range("AN4").Value = 6
$ws.Range("AO4").Value = 26
$ws.Range("AP4").Value = 41
$ws.Range("AQ4").Value = 101
$ws.Range("AR4").Value = 151
$ws.Range("AS4").Value = 401
$ws.Range("AT4").Value = 2.25
$ws.Range("AW4").Value = 3.75
$ws.Range("AX4").Value = 11
$ws.Range("AY4").Value = 26
$ws.Range("AZ4").Value = 41
$ws.Range("BA4").Value = 67
$ws.Range("BB4").Value = 251
$ws.Range("BC4").Value = 126
$ws.Range("BD4").Value = 126
$ws.Range("N5").Value = 8
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.57
$ws.Range("AM5").Value = 900
